$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (cluster-pair rows), matching the refreshed Il12b-Il12rb2 NATMI output
# Columns: A..T => Sending cluster, Ligand symbol, Receptor symbol, Target cluster, then 16 numeric metrics
$rowsData = @(
    @("FAPs", "Il12b", "Il12rb2", "ECs", 1, 0.3333333333333333, 0.091959, 0.275877, 0.06282071922171414, 0.06282071922171413, 3, 1, 0.3922993333333333, 1.176898, 0.1197958591217032, 0.1197958591217032, 0.036075454394, 0.324679089546, 0.007525662029808538, 0.007525662029808539),
    @("FAPs", "Il12b", "Il12rb2", "FAPs", 1, 0.3333333333333333, 0.091959, 0.275877, 0.06282071922171414, 0.06282071922171413, 3, 1, 0.544463, 1.633389, 0.1662618498246574, 0.1662618498246574, 0.050068273017, 0.450614457153, 0.01044468898511761, 0.0104446889851176),
    @("FAPs", "Il12b", "Il12rb2", "M2", 1, 0.3333333333333333, 0.091959, 0.275877, 0.06282071922171414, 0.06282071922171413, 3, 1, 1.882823333333333, 5.64847, 0.5749549377882933, 0.5749549377882933, 0.17314255091, 1.55828295819, 0.0361190827119365, 0.03611908271193649),
    @("FAPs", "Il12b", "Il12rb2", "sCs", 1, 0.3333333333333333, 0.091959, 0.275877, 0.06282071922171414, 0.06282071922171413, 3, 1, 0.4551463333333333, 1.365439, 0.1389873532653461, 0.1389873532653461, 0.04185480166699999, 0.376693215003, 0.008731285494851498, 0.008731285494851498),
    @("M2", "Il12b", "Il12rb2", "ECs", 3, 1, 1.371873333333334, 4.115620000000001, 0.9371792807782859, 0.9371792807782858, 3, 1, 0.3922993333333333, 1.176898, 0.1197958591217032, 0.1197958591217032, 0.5381849940844445, 4.843664946760001, 0.1122701970918946, 0.1122701970918947),
    @("M2", "Il12b", "Il12rb2", "FAPs", 3, 1, 1.371873333333334, 4.115620000000001, 0.9371792807782859, 0.9371792807782858, 3, 1, 0.544463, 1.633389, 0.1662618498246574, 0.1662618498246574, 0.7469342706866668, 6.722408436180001, 0.1558171608395398, 0.1558171608395398),
    @("M2", "Il12b", "Il12rb2", "M2", 3, 1, 1.371873333333334, 4.115620000000001, 0.9371792807782859, 0.9371792807782858, 3, 1, 1.882823333333333, 5.64847, 0.5749549377882933, 0.5749549377882933, 2.582995122377778, 23.2469561014, 0.5388358550763569, 0.5388358550763568),
    @("M2", "Il12b", "Il12rb2", "sCs", 3, 1, 1.371873333333334, 4.115620000000001, 0.9371792807782859, 0.9371792807782858, 3, 1, 0.4551463333333333, 1.365439, 0.1389873532653461, 0.1389873532653461, 0.6244031174644444, 5.61962805718, 0.1302560677704946, 0.1302560677704946)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rowsData[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
